$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.178.52'
Set-TextValue 'E2' '  -2.53%  '
Set-TextValue 'D3' '2.380.57'
Set-TextValue 'E3' '  -2.59%  '
Set-TextValue 'E4' '  +0.32%  '
Set-TextValue 'D5' '562.31'
Set-TextValue 'E5' '  -2.71%  '
Set-TextValue 'D6' '137.54'
Set-TextValue 'E6' '  -2.19%  '
Set-TextValue 'E7' '  +0.01%  '
Set-TextValue 'D8' '0.536'
Set-TextValue 'E8' '  +0.60%  '
Set-TextValue 'D9' '2.379.01'
Set-TextValue 'E9' '  -2.26%  '
Set-TextValue 'D10' '0.104'
Set-TextValue 'E10' '  -4.50%  '
Set-TextValue 'D11' '0.159'
Set-TextValue 'E11' '  -1.07%  '
Set-TextValue 'D12' '5.06'
Set-TextValue 'E12' '  -1.71%  '
Set-TextValue 'D13' '0.336'
Set-TextValue 'E13' '  -1.07%  '
Set-TextValue 'D14' '25.54'
Set-TextValue 'E14' '  -1.09%  '
Set-TextValue 'D15' '2.817.15'
Set-TextValue 'E15' '  -2.53%  '
Set-TextValue 'D16' '0.0000165'
Set-TextValue 'E16' '  -3.52%  '
Set-TextValue 'D17' '60.366.12'
Set-TextValue 'E17' '  -2.17%  '
Set-TextValue 'D18' '2.388.74'
Set-TextValue 'E18' '  -2.19%  '
Set-TextValue 'D19' '8.28'
Set-TextValue 'E19' '  +14.63%  '
Set-TextValue 'D20' '10.54'
Set-TextValue 'E20' '  -0.27%  '
Set-TextValue 'D21' '324.02'
Set-TextValue 'E21' '  -0.23%  '
Set-TextValue 'D22' '4.02'
Set-TextValue 'E22' '  -1.18%  '
Set-TextValue 'D23' '6.10'
Set-TextValue 'E23' '  +0.81%  '
Set-TextValue 'E24' '  -0.05%  '
Set-TextValue 'D25' '1.77'
Set-TextValue 'E25' '  -8.28%  '
Set-TextValue 'D26' '64.35'
Set-TextValue 'E26' '  -0.80%  '
Set-TextValue 'D27' '549.86'
Set-TextValue 'E27' '  -5.04%  '
Set-TextValue 'D28' '7.98'
Set-TextValue 'E28' '  -12.37%  '
Set-TextValue 'D29' '2.522.52'
Set-TextValue 'E29' '  -1.77%  '
Set-TextValue 'D30' '0.0₃0894'
Set-TextValue 'E30' '  -2.60%  '
Set-TextValue 'D31' '7.85'
Set-TextValue 'E31' '  -0.78%  '
Set-TextValue 'D32' '1.28'
Set-TextValue 'E32' '  -4.86%  '
Set-TextValue 'D33' '1.79'
Set-TextValue 'E33' '  -3.99%  '
Set-TextValue 'D35' '1.00'
Set-TextValue 'E35' '  -0.22%  '
Set-TextValue 'D36' '153.98'
Set-TextValue 'E36' '  +1.55%  '
Set-TextValue 'D37' '1.40'
Set-TextValue 'E37' '  +0.70%  '
Set-TextValue 'D38' '0.365'
Set-TextValue 'E38' '  -1.66%  '
Set-TextValue 'D39' '4.49'
Set-TextValue 'E39' '  -4.28%  '
Set-TextValue 'D40' '18.20'
Set-TextValue 'E40' '  -0.36%  '
Set-TextValue 'D41' '5.01'
Set-TextValue 'E41' '  -2.10%  '
Set-TextValue 'D43' '1.62'
Set-TextValue 'E43' '  -3.52%  '
Set-TextValue 'D44' '2.27'
Set-TextValue 'E44' '  -3.58%  '
Set-TextValue 'D45' '0.0₆0278'
Set-TextValue 'E45' '  -2.97%  '
Set-TextValue 'D46' '142.38'
Set-TextValue 'E46' '  +0.08%  '
Set-TextValue 'D47' '3.47'
Set-TextValue 'E47' '  -2.50%  '
Set-TextValue 'D48' '0.582'
Set-TextValue 'E48' '  -2.52%  '
Set-TextValue 'D49' '0.0494'
Set-TextValue 'E49' '  -2.77%  '
Set-TextValue 'D50' '18.77'
Set-TextValue 'E50' '  -3.90%  '
Set-TextValue 'D51' '0.0894'
Set-TextValue 'E51' '  -0.23%  '
